$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '30.121.13'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').Value = '1.643.85'
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'215.94"
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = "'29.40"
$ws.Range('E8').Value = '  +5.66%  '
$ws.Range('E9').Value = '  +4.04%  '
$ws.Range('D10').Value = "'0.0616"
$ws.Range('E10').Value = '  +2.30%  '
$ws.Range('D11').Value = "'0.0916"
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').Value = '1.876.93'
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').Value = '1.641.12'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').Value = "'0.576"
$ws.Range('E14').Value = '  +5.81%  '
$ws.Range('D15').Value = "'9.51"
$ws.Range('E15').Value = '  +23.38%  '
$ws.Range('D16').Value = "'3.93"
$ws.Range('E16').Value = '  +4.85%  '
$ws.Range('D17').Value = '30.123.17'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('D18').Value = "'65.11"
$ws.Range('E18').Value = '  +2.09%  '
$ws.Range('D19').Value = "'248.36"
$ws.Range('E19').Value = '  +2.85%  '
$ws.Range('E20').Value = '  +2.34%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = "'4.21"
$ws.Range('E22').Value = '  +5.15%  '
$ws.Range('D23').Value = "'9.93"
$ws.Range('E23').Value = '  +6.62%  '
$ws.Range('D24').Value = "'2.15"
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('D25').Value = "'159.12"
$ws.Range('D26').Value = "'15.80"
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('E27').Value = '  +2.87%  '
$ws.Range('E28').Value = '  +4.23%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = "'0.0494"
$ws.Range('E30').Value = '  +2.83%  '
$ws.Range('E31').Value = '  +6.26%  '
$ws.Range('D32').Value = "'3.44"
$ws.Range('E32').Value = '  +6.34%  '
$ws.Range('D33').Value = "'3.22"
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').Value = '1.440.28'
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('E35').Value = '  +7.50%  '
$ws.Range('E36').Value = '  +2.03%  '
$ws.Range('D37').Value = "'2.87"
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').Value = "'77.97"
$ws.Range('E38').Value = '  +18.01%  '
$ws.Range('E39').Value = '  +2.16%  '
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('D41').Value = "'0.560"
$ws.Range('E41').Value = '  +2.59%  '
$ws.Range('E42').Value = '  +3.23%  '
$ws.Range('D43').Value = "'0.847"
$ws.Range('E43').Value = '  +3.93%  '
$ws.Range('D44').Value = "'55.84"
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('E46').Value = '  +5.24%  '
$ws.Range('D47').Value = "'0.999"
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = "'5.40"
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').Value = '1.783.08'
$ws.Range('E49').Value = '  +2.38%  '
$ws.Range('D50').Value = '0.0₆0115'
$ws.Range('E50').Value = '  +11.28%  '
$ws.Range('D51').Value = "'90.55"
$ws.Range('E51').Value = '  +4.51%  '

# The quote-prefix trick above marks the cell with a 'quotePrefix' style so Excel
# stops treating it as a number. Reset style back to Normal (no style index) to
# match the original plain/unstyled data cells once the text value is locked in.
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D51').Style = "Normal"
